# Update the "Speedup Results" benchmark table (rows 2-31, columns A-J)
# with refreshed benchmark numbers (adds spherical timings / regroups by
# covariance type per the updated benchmark run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: _maximization_step / diag
$ws.Cells.Item(2, 1).Value = "_maximization_step"
$ws.Cells.Item(2, 2).Value = "diag"
$ws.Cells.Item(2, 3).Value = 500
$ws.Cells.Item(2, 4).Value = 20
$ws.Cells.Item(2, 5).Value = 5
$ws.Cells.Item(2, 6).Value = 0.450134911807254
$ws.Cells.Item(2, 7).Value = 0.1324756908531867
$ws.Cells.Item(2, 8).Value = 0.1829224987886846
$ws.Cells.Item(2, 9).Value = 0.0206478852027326
$ws.Cells.Item(2, 10).Value = 2.460795773007988

# Row 3: TorchGaussianMixture.fit / diag
$ws.Cells.Item(3, 1).Value = "TorchGaussianMixture.fit"
$ws.Cells.Item(3, 2).Value = "diag"
$ws.Cells.Item(3, 3).Value = 500
$ws.Cells.Item(3, 4).Value = 20
$ws.Cells.Item(3, 5).Value = 5
$ws.Cells.Item(3, 6).Value = 13.52092534458886
$ws.Cells.Item(3, 7).Value = 1.191474979905249
$ws.Cells.Item(3, 8).Value = 11.64670833774532
$ws.Cells.Item(3, 9).Value = 2.490389354718123
$ws.Cells.Item(3, 10).Value = 1.160922464312897

# Row 4: _maximization_step / diag
$ws.Cells.Item(4, 1).Value = "_maximization_step"
$ws.Cells.Item(4, 2).Value = "diag"
$ws.Cells.Item(4, 3).Value = 1000
$ws.Cells.Item(4, 4).Value = 50
$ws.Cells.Item(4, 5).Value = 10
$ws.Cells.Item(4, 6).Value = 2.331453916849568
$ws.Cells.Item(4, 7).Value = 1.984161536843239
$ws.Cells.Item(4, 8).Value = 1.353842503158376
$ws.Cells.Item(4, 9).Value = 0.3296520629571567
$ws.Cells.Item(4, 10).Value = 1.722101286826589

# Row 5: TorchGaussianMixture.fit / diag
$ws.Cells.Item(5, 1).Value = "TorchGaussianMixture.fit"
$ws.Cells.Item(5, 2).Value = "diag"
$ws.Cells.Item(5, 3).Value = 1000
$ws.Cells.Item(5, 4).Value = 50
$ws.Cells.Item(5, 5).Value = 5
$ws.Cells.Item(5, 6).Value = 37.5325593049638
$ws.Cells.Item(5, 7).Value = 6.683749062159406
$ws.Cells.Item(5, 8).Value = 42.39961568964645
$ws.Cells.Item(5, 9).Value = 12.22135834763968
$ws.Cells.Item(5, 10).Value = 0.8852098938747895

# Row 6: _maximization_step / diag
$ws.Cells.Item(6, 1).Value = "_maximization_step"
$ws.Cells.Item(6, 2).Value = "diag"
$ws.Cells.Item(6, 3).Value = 2000
$ws.Cells.Item(6, 4).Value = 100
$ws.Cells.Item(6, 5).Value = 20
$ws.Cells.Item(6, 6).Value = 11.01770000532269
$ws.Cells.Item(6, 7).Value = 1.511694376542218
$ws.Cells.Item(6, 8).Value = 27.77586799347773
$ws.Cells.Item(6, 9).Value = 2.239002524270575
$ws.Cells.Item(6, 10).Value = 0.3966644717605169

# Row 7: _estimate_log_gaussian_prob_tied_precchol / tied
$ws.Cells.Item(7, 1).Value = "_estimate_log_gaussian_prob_tied_precchol"
$ws.Cells.Item(7, 2).Value = "tied"
$ws.Cells.Item(7, 3).Value = 500
$ws.Cells.Item(7, 4).Value = 20
$ws.Cells.Item(7, 5).Value = 5
$ws.Cells.Item(7, 6).Value = 0.4193803993985057
$ws.Cells.Item(7, 7).Value = 0.2386536237897989
$ws.Cells.Item(7, 8).Value = 0.2107206964865327
$ws.Cells.Item(7, 9).Value = 0.05275635152072836
$ws.Cells.Item(7, 10).Value = 1.990219311112179

# Row 8: _maximization_step / tied
$ws.Cells.Item(8, 1).Value = "_maximization_step"
$ws.Cells.Item(8, 2).Value = "tied"
$ws.Cells.Item(8, 3).Value = 500
$ws.Cells.Item(8, 4).Value = 20
$ws.Cells.Item(8, 5).Value = 5
$ws.Cells.Item(8, 6).Value = 0.3138965112157166
$ws.Cells.Item(8, 7).Value = 0.04981436749508657
$ws.Cells.Item(8, 8).Value = 0.3371396102011204
$ws.Cells.Item(8, 9).Value = 0.04630397406384053
$ws.Cells.Item(8, 10).Value = 0.9310579407399263

# Row 9: TorchGaussianMixture.fit / tied
$ws.Cells.Item(9, 1).Value = "TorchGaussianMixture.fit"
$ws.Cells.Item(9, 2).Value = "tied"
$ws.Cells.Item(9, 3).Value = 500
$ws.Cells.Item(9, 4).Value = 20
$ws.Cells.Item(9, 5).Value = 5
$ws.Cells.Item(9, 6).Value = 4.949284659232944
$ws.Cells.Item(9, 7).Value = 1.965017131183696
$ws.Cells.Item(9, 8).Value = 2.894892008043826
$ws.Cells.Item(9, 9).Value = 0.4750682050116928
$ws.Cells.Item(9, 10).Value = 1.709661239687259

# Row 10: _estimate_log_gaussian_prob_tied_precchol / tied
$ws.Cells.Item(10, 1).Value = "_estimate_log_gaussian_prob_tied_precchol"
$ws.Cells.Item(10, 2).Value = "tied"
$ws.Cells.Item(10, 3).Value = 1000
$ws.Cells.Item(10, 4).Value = 50
$ws.Cells.Item(10, 5).Value = 10
$ws.Cells.Item(10, 6).Value = 2.016141905914992
$ws.Cells.Item(10, 7).Value = 0.4930243904024206
$ws.Cells.Item(10, 8).Value = 2.465680608293042
$ws.Cells.Item(10, 9).Value = 0.9145870598931194
$ws.Cells.Item(10, 10).Value = 0.8176816977567669

# Row 11: _maximization_step / tied
$ws.Cells.Item(11, 1).Value = "_maximization_step"
$ws.Cells.Item(11, 2).Value = "tied"
$ws.Cells.Item(11, 3).Value = 1000
$ws.Cells.Item(11, 4).Value = 50
$ws.Cells.Item(11, 5).Value = 10
$ws.Cells.Item(11, 6).Value = 3.175517916679382
$ws.Cells.Item(11, 7).Value = 1.316052945837307
$ws.Cells.Item(11, 8).Value = 3.607950697187334
$ws.Cells.Item(11, 9).Value = 0.9826375094639167
$ws.Cells.Item(11, 10).Value = 0.8801444873276496

# Row 12: TorchGaussianMixture.fit / tied
$ws.Cells.Item(12, 1).Value = "TorchGaussianMixture.fit"
$ws.Cells.Item(12, 2).Value = "tied"
$ws.Cells.Item(12, 3).Value = 1000
$ws.Cells.Item(12, 4).Value = 50
$ws.Cells.Item(12, 5).Value = 5
$ws.Cells.Item(12, 6).Value = 11.40718901297078
$ws.Cells.Item(12, 7).Value = 2.895266571640174
$ws.Cells.Item(12, 8).Value = 9.927896996183941
$ws.Cells.Item(12, 9).Value = 2.245066842066469
$ws.Cells.Item(12, 10).Value = 1.149003562119494

# Row 13: _estimate_log_gaussian_prob_tied_precchol / tied
$ws.Cells.Item(13, 1).Value = "_estimate_log_gaussian_prob_tied_precchol"
$ws.Cells.Item(13, 2).Value = "tied"
$ws.Cells.Item(13, 3).Value = 2000
$ws.Cells.Item(13, 4).Value = 100
$ws.Cells.Item(13, 5).Value = 20
$ws.Cells.Item(13, 6).Value = 15.03023460390978
$ws.Cells.Item(13, 7).Value = 1.309180928019012
$ws.Cells.Item(13, 8).Value = 25.91667589731514
$ws.Cells.Item(13, 9).Value = 5.240738662265411
$ws.Cells.Item(13, 10).Value = 0.5799445370024032

# Row 14: _maximization_step / tied
$ws.Cells.Item(14, 1).Value = "_maximization_step"
$ws.Cells.Item(14, 2).Value = "tied"
$ws.Cells.Item(14, 3).Value = 2000
$ws.Cells.Item(14, 4).Value = 100
$ws.Cells.Item(14, 5).Value = 20
$ws.Cells.Item(14, 6).Value = 25.57757350150496
$ws.Cells.Item(14, 7).Value = 1.886056339373293
$ws.Cells.Item(14, 8).Value = 80.78049459727481
$ws.Cells.Item(14, 9).Value = 3.858457963010033
$ws.Cells.Item(14, 10).Value = 0.3166305632197483

# Row 15: _compute_precisions_cholesky / full
$ws.Cells.Item(15, 1).Value = "_compute_precisions_cholesky"
$ws.Cells.Item(15, 2).Value = "full"
$ws.Cells.Item(15, 3).Value = 500
$ws.Cells.Item(15, 4).Value = 20
$ws.Cells.Item(15, 5).Value = 5
$ws.Cells.Item(15, 6).Value = 0.1666151045355946
$ws.Cells.Item(15, 7).Value = 0.1135115570616916
$ws.Cells.Item(15, 8).Value = 0.1031568099278957
$ws.Cells.Item(15, 9).Value = 0.02668522680127083
$ws.Cells.Item(15, 10).Value = 1.615163406585129

# Row 16: _compute_precisions / full
$ws.Cells.Item(16, 1).Value = "_compute_precisions"
$ws.Cells.Item(16, 2).Value = "full"
$ws.Cells.Item(16, 3).Value = 500
$ws.Cells.Item(16, 4).Value = 20
$ws.Cells.Item(16, 5).Value = 5
$ws.Cells.Item(16, 6).Value = 0.03105999203398824
$ws.Cells.Item(16, 7).Value = 0.000678339726004085
$ws.Cells.Item(16, 8).Value = 0.007837096927687526
$ws.Cells.Item(16, 9).Value = 0.0002569159501781693
$ws.Cells.Item(16, 10).Value = 3.96320121093016

# Row 17: _estimate_log_gaussian_prob_full_precchol / full
$ws.Cells.Item(17, 1).Value = "_estimate_log_gaussian_prob_full_precchol"
$ws.Cells.Item(17, 2).Value = "full"
$ws.Cells.Item(17, 3).Value = 500
$ws.Cells.Item(17, 4).Value = 20
$ws.Cells.Item(17, 5).Value = 5
$ws.Cells.Item(17, 6).Value = 0.4364194988738745
$ws.Cells.Item(17, 7).Value = 0.06113912044385759
$ws.Cells.Item(17, 8).Value = 0.2589675015769899
$ws.Cells.Item(17, 9).Value = 0.02639495297570719
$ws.Cells.Item(17, 10).Value = 1.685228826846171

# Row 18: _maximization_step / full
$ws.Cells.Item(18, 1).Value = "_maximization_step"
$ws.Cells.Item(18, 2).Value = "full"
$ws.Cells.Item(18, 3).Value = 500
$ws.Cells.Item(18, 4).Value = 20
$ws.Cells.Item(18, 5).Value = 5
$ws.Cells.Item(18, 6).Value = 1.848321396391839
$ws.Cells.Item(18, 7).Value = 0.8135488201834762
$ws.Cells.Item(18, 8).Value = 1.162103796377778
$ws.Cells.Item(18, 9).Value = 0.8319351821421227
$ws.Cells.Item(18, 10).Value = 1.590495962712598

# Row 19: TorchGaussianMixture.fit / full
$ws.Cells.Item(19, 1).Value = "TorchGaussianMixture.fit"
$ws.Cells.Item(19, 2).Value = "full"
$ws.Cells.Item(19, 3).Value = 500
$ws.Cells.Item(19, 4).Value = 20
$ws.Cells.Item(19, 5).Value = 5
$ws.Cells.Item(19, 6).Value = 29.92165767742942
$ws.Cells.Item(19, 7).Value = 4.04296688410544
$ws.Cells.Item(19, 8).Value = 21.52464932684476
$ws.Cells.Item(19, 9).Value = 3.778393908502379
$ws.Cells.Item(19, 10).Value = 1.390111273037662

# Row 20: _compute_precisions_cholesky / full
$ws.Cells.Item(20, 1).Value = "_compute_precisions_cholesky"
$ws.Cells.Item(20, 2).Value = "full"
$ws.Cells.Item(20, 3).Value = 1000
$ws.Cells.Item(20, 4).Value = 50
$ws.Cells.Item(20, 5).Value = 10
$ws.Cells.Item(20, 6).Value = 0.5656105000525713
$ws.Cells.Item(20, 7).Value = 0.3964756663224189
$ws.Cells.Item(20, 8).Value = 0.1705931150354445
$ws.Cells.Item(20, 9).Value = 0.048250049541323
$ws.Cells.Item(20, 10).Value = 3.31555291627715

# Row 21: _compute_precisions / full
$ws.Cells.Item(21, 1).Value = "_compute_precisions"
$ws.Cells.Item(21, 2).Value = "full"
$ws.Cells.Item(21, 3).Value = 1000
$ws.Cells.Item(21, 4).Value = 50
$ws.Cells.Item(21, 5).Value = 10
$ws.Cells.Item(21, 6).Value = 0.1453039061743766
$ws.Cells.Item(21, 7).Value = 0.02700967109957315
$ws.Cells.Item(21, 8).Value = 0.04763819160871208
$ws.Cells.Item(21, 9).Value = 0.01186937266446013
$ws.Cells.Item(21, 10).Value = 3.05015579449081

# Row 22: _estimate_log_gaussian_prob_full_precchol / full
$ws.Cells.Item(22, 1).Value = "_estimate_log_gaussian_prob_full_precchol"
$ws.Cells.Item(22, 2).Value = "full"
$ws.Cells.Item(22, 3).Value = 1000
$ws.Cells.Item(22, 4).Value = 50
$ws.Cells.Item(22, 5).Value = 10
$ws.Cells.Item(22, 6).Value = 2.888031402835622
$ws.Cells.Item(22, 7).Value = 0.7224816116141111
$ws.Cells.Item(22, 8).Value = 2.042149292537943
$ws.Cells.Item(22, 9).Value = 0.4050252242186253
$ws.Cells.Item(22, 10).Value = 1.414211690295392

# Row 23: _maximization_step / full
$ws.Cells.Item(23, 1).Value = "_maximization_step"
$ws.Cells.Item(23, 2).Value = "full"
$ws.Cells.Item(23, 3).Value = 1000
$ws.Cells.Item(23, 4).Value = 50
$ws.Cells.Item(23, 5).Value = 10
$ws.Cells.Item(23, 6).Value = 6.503827008418739
$ws.Cells.Item(23, 7).Value = 2.073988007254053
$ws.Cells.Item(23, 8).Value = 3.233443002682179
$ws.Cells.Item(23, 9).Value = 1.290858187326016
$ws.Cells.Item(23, 10).Value = 2.011424664985199

# Row 24: TorchGaussianMixture.fit / full
$ws.Cells.Item(24, 1).Value = "TorchGaussianMixture.fit"
$ws.Cells.Item(24, 2).Value = "full"
$ws.Cells.Item(24, 3).Value = 1000
$ws.Cells.Item(24, 4).Value = 50
$ws.Cells.Item(24, 5).Value = 5
$ws.Cells.Item(24, 6).Value = 119.786839990411
$ws.Cells.Item(24, 7).Value = 27.10058913911693
$ws.Cells.Item(24, 8).Value = 61.60724300813551
$ws.Cells.Item(24, 9).Value = 4.201825111264215
$ws.Cells.Item(24, 10).Value = 1.94436293756227

# Row 25: _compute_precisions_cholesky / full
$ws.Cells.Item(25, 1).Value = "_compute_precisions_cholesky"
$ws.Cells.Item(25, 2).Value = "full"
$ws.Cells.Item(25, 3).Value = 2000
$ws.Cells.Item(25, 4).Value = 100
$ws.Cells.Item(25, 5).Value = 5
$ws.Cells.Item(25, 6).Value = 0.4165725084021688
$ws.Cells.Item(25, 7).Value = 0.03748828564352247
$ws.Cells.Item(25, 8).Value = 0.3259302116930485
$ws.Cells.Item(25, 9).Value = 0.09037528925134744
$ws.Cells.Item(25, 10).Value = 1.278103389797091

# Row 26: _compute_precisions / full
$ws.Cells.Item(26, 1).Value = "_compute_precisions"
$ws.Cells.Item(26, 2).Value = "full"
$ws.Cells.Item(26, 3).Value = 2000
$ws.Cells.Item(26, 4).Value = 100
$ws.Cells.Item(26, 5).Value = 5
$ws.Cells.Item(26, 6).Value = 0.2027027891017497
$ws.Cells.Item(26, 7).Value = 0.03430058066310663
$ws.Cells.Item(26, 8).Value = 0.1254188013263047
$ws.Cells.Item(26, 9).Value = 0.009806133452759153
$ws.Cells.Item(26, 10).Value = 1.616207354544664

# Row 27: _estimate_log_gaussian_prob_full_precchol / full
$ws.Cells.Item(27, 1).Value = "_estimate_log_gaussian_prob_full_precchol"
$ws.Cells.Item(27, 2).Value = "full"
$ws.Cells.Item(27, 3).Value = 2000
$ws.Cells.Item(27, 4).Value = 100
$ws.Cells.Item(27, 5).Value = 5
$ws.Cells.Item(27, 6).Value = 5.245532694971189
$ws.Cells.Item(27, 7).Value = 0.7317313876458609
$ws.Cells.Item(27, 8).Value = 5.775378999533132
$ws.Cells.Item(27, 9).Value = 0.3615964659091776
$ws.Cells.Item(27, 10).Value = 0.9082577429801967

# Row 28: _maximization_step / full
$ws.Cells.Item(28, 1).Value = "_maximization_step"
$ws.Cells.Item(28, 2).Value = "full"
$ws.Cells.Item(28, 3).Value = 2000
$ws.Cells.Item(28, 4).Value = 100
$ws.Cells.Item(28, 5).Value = 5
$ws.Cells.Item(28, 6).Value = 10.50152819952928
$ws.Cells.Item(28, 7).Value = 1.650155901872884
$ws.Cells.Item(28, 8).Value = 8.407217304920778
$ws.Cells.Item(28, 9).Value = 1.403609561785657
$ws.Cells.Item(28, 10).Value = 1.249108690622602

# Row 29: _kmeans_lloyd_with_init / N/A
$ws.Cells.Item(29, 1).Value = "_kmeans_lloyd_with_init"
$ws.Cells.Item(29, 2).Value = "N/A"
$ws.Cells.Item(29, 3).Value = 500
$ws.Cells.Item(29, 4).Value = 20
$ws.Cells.Item(29, 5).Value = 5
$ws.Cells.Item(29, 6).Value = 4.854389210231602
$ws.Cells.Item(29, 7).Value = 3.393013982703177
$ws.Cells.Item(29, 8).Value = 2.051664609462023
$ws.Cells.Item(29, 9).Value = 0.1651375077439367
$ws.Cells.Item(29, 10).Value = 2.366073474116462

# Row 30: _kmeans_lloyd_with_init / N/A
$ws.Cells.Item(30, 1).Value = "_kmeans_lloyd_with_init"
$ws.Cells.Item(30, 2).Value = "N/A"
$ws.Cells.Item(30, 3).Value = 1000
$ws.Cells.Item(30, 4).Value = 50
$ws.Cells.Item(30, 5).Value = 10
$ws.Cells.Item(30, 6).Value = 9.482237393967807
$ws.Cells.Item(30, 7).Value = 2.0486794635028
$ws.Cells.Item(30, 8).Value = 7.57343020522967
$ws.Cells.Item(30, 9).Value = 1.344112971547114
$ws.Cells.Item(30, 10).Value = 1.252039978848693

# Row 31: _kmeans_lloyd_with_init / N/A
$ws.Cells.Item(31, 1).Value = "_kmeans_lloyd_with_init"
$ws.Cells.Item(31, 2).Value = "N/A"
$ws.Cells.Item(31, 3).Value = 2000
$ws.Cells.Item(31, 4).Value = 100
$ws.Cells.Item(31, 5).Value = 20
$ws.Cells.Item(31, 6).Value = 27.11349718738347
$ws.Cells.Item(31, 7).Value = 6.044830419891486
$ws.Cells.Item(31, 8).Value = 13.96285960217938
$ws.Cells.Item(31, 9).Value = 2.868270668984783
$ws.Cells.Item(31, 10).Value = 1.941829822821644

